$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the style of the existing
# header cells (e.g. G1) so it matches the rest of the header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the new data value in H2 (plain numeric cell, no special style).
$ws.Range("H2").Value = 1
